# "added test de perception" — fill in Inessa Kechek's answers to the
# Belbin team-role self-perception questionnaire and fix a typo.

$wb = $excel.ActiveWorkbook
$wsQ = $wb.Worksheets.Item("Questionnaire")
$wsR = $wb.Worksheets.Item("Résultats")

# --- Fix a spelling typo in one of the questionnaire prompts --------------
# "conplexes" -> "complexes"
$wsQ.Range("B78").Value = "Je n'arrive pas toujours à expliquer certains aspects complexes d'un sujet"

# --- Fill in the answers (point values 1-5) entered for each question -----
$answers = @{
    "C8"  = 3
    "C10" = 3
    "C11" = 4
    "C19" = 4
    "C20" = 4
    "C21" = 2
    "C30" = 2
    "C35" = 5
    "C37" = 3
    "C41" = 5
    "C46" = 3
    "C48" = 2
    "C52" = 3
    "C53" = 3
    "C57" = 4
    "C66" = 4
    "C68" = 4
    "C70" = 2
    "C75" = 4
    "C76" = 3
    "C79" = 3
}

foreach ($addr in $answers.Keys) {
    $wsQ.Range($addr).Value = $answers[$addr]
}

# --- Restore/update the view state (zoom + selection + active sheet) ------
$wsQ.Activate()
$excel.ActiveWindow.Zoom = 87
$wsQ.Range("H83").Select()

$wsR.Activate()
$wsR.Range("B36:B37").Select()
